$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18,8).Value = 2567
$ws.Cells.Item(18,9).Value = 2567
$ws.Cells.Item(18,10).Value = 0
$ws.Cells.Item(18,11).Value = 2567
$ws.Cells.Item(18,12).Value = 0
$ws.Cells.Item(18,13).Value = -2283

$ws.Cells.Item(19,8).Value = 1909.4
$ws.Cells.Item(19,9).Value = 2499
$ws.Cells.Item(19,10).Value = 1762
$ws.Cells.Item(19,11).Value = 2499
$ws.Cells.Item(19,12).Value = 1762
$ws.Cells.Item(19,13).Value = -2324

$ws.Cells.Item(40,8).Value = 819.6
$ws.Cells.Item(40,9).Value = 532.6667
$ws.Cells.Item(40,10).Value = 1250
$ws.Cells.Item(40,11).Value = 532.6667
$ws.Cells.Item(40,12).Value = 1250
$ws.Cells.Item(40,13).Value = -357.6667
$ws.Cells.Item(40,14).Value = -1600

$ws.Cells.Item(53,8).Value = 113.5
$ws.Cells.Item(53,9).Value = 96.14286
$ws.Cells.Item(53,10).Value = 130.85715
$ws.Cells.Item(53,11).Value = 96.14286
$ws.Cells.Item(53,12).Value = 130.85715
$ws.Cells.Item(53,13).Value = 540.85714

$ws.Cells.Item(64,8).Value = 4728.4287
$ws.Cells.Item(64,9).Value = 4456.857
$ws.Cells.Item(64,10).Value = 5000
$ws.Cells.Item(64,11).Value = 4456.857
$ws.Cells.Item(64,12).Value = 5000
$ws.Cells.Item(64,13).Value = -4208.857
$ws.Cells.Item(64,14).Value = -5496

$ws.Cells.Item(67,8).Value = 4728.4287
$ws.Cells.Item(67,9).Value = 4456.857
$ws.Cells.Item(67,10).Value = 5000
$ws.Cells.Item(67,11).Value = 4456.857
$ws.Cells.Item(67,12).Value = 5000
$ws.Cells.Item(67,13).Value = -3598.857
$ws.Cells.Item(67,14).Value = -6716

$ws.Cells.Item(74,8).Value = 3171.2222
$ws.Cells.Item(74,9).Value = 3171.2222
$ws.Cells.Item(74,10).Value = 0
$ws.Cells.Item(74,11).Value = 3171.2222
$ws.Cells.Item(74,12).Value = 0
$ws.Cells.Item(74,13).Value = -2235.2222

$ws.Cells.Item(77,8).Value = 3171.2222
$ws.Cells.Item(77,9).Value = 3171.2222
$ws.Cells.Item(77,10).Value = 0
$ws.Cells.Item(77,11).Value = 15856.111
$ws.Cells.Item(77,12).Value = 0
$ws.Cells.Item(77,13).Value = -11176.111

$ws.Cells.Item(113,8).Value = 9732.625
$ws.Cells.Item(113,9).Value = 6287.6665
$ws.Cells.Item(113,10).Value = 11799.6
$ws.Cells.Item(113,11).Value = 6287.6665
$ws.Cells.Item(113,12).Value = 11799.6
$ws.Cells.Item(113,13).Value = -3033.6665
$ws.Cells.Item(113,14).Value = -18307.6

$ws.Cells.Item(115,8).Value = 377
$ws.Cells.Item(115,9).Value = 377
$ws.Cells.Item(115,10).Value = 0
$ws.Cells.Item(115,11).Value = 1131
$ws.Cells.Item(115,12).Value = 0
$ws.Cells.Item(115,13).Value = 436

$ws.Cells.Item(132,8).Value = 1980.5
$ws.Cells.Item(132,9).Value = 2180.3333
$ws.Cells.Item(132,10).Value = 581.6667
$ws.Cells.Item(132,11).Value = 6540.999899999999
$ws.Cells.Item(132,12).Value = 1745.0001
$ws.Cells.Item(132,13).Value = -4010.999899999999
$ws.Cells.Item(132,14).Value = -6805.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32,8).Value = 1453.0952
$ws.Cells.Item(32,9).Value = 1280.6102
$ws.Cells.Item(32,10).Value = 3997.25
$ws.Cells.Item(32,11).Value = 1280.6102
$ws.Cells.Item(32,12).Value = 3997.25
$ws.Cells.Item(32,13).Value = -993.6102000000001
$ws.Cells.Item(32,14).Value = -4571.25

$ws.Cells.Item(61,8).Value = 3926.4
$ws.Cells.Item(61,9).Value = 3926.4
$ws.Cells.Item(61,10).Value = 0
$ws.Cells.Item(61,11).Value = 3926.4
$ws.Cells.Item(61,12).Value = 0
$ws.Cells.Item(61,13).Value = -3714.4

$ws.Cells.Item(63,8).Value = 3849.125
$ws.Cells.Item(63,9).Value = 3623.25
$ws.Cells.Item(63,10).Value = 4075
$ws.Cells.Item(63,11).Value = 3623.25
$ws.Cells.Item(63,12).Value = 4075
$ws.Cells.Item(63,13).Value = -2937.25
$ws.Cells.Item(63,14).Value = -5447

$ws.Cells.Item(66,8).Value = 3849.125
$ws.Cells.Item(66,9).Value = 3623.25
$ws.Cells.Item(66,10).Value = 4075
$ws.Cells.Item(66,11).Value = 18116.25
$ws.Cells.Item(66,12).Value = 20375
$ws.Cells.Item(66,13).Value = -14684.25
$ws.Cells.Item(66,14).Value = -27239

$ws.Cells.Item(102,8).Value = 1099.5714
$ws.Cells.Item(102,9).Value = 917.6667
$ws.Cells.Item(102,10).Value = 2191
$ws.Cells.Item(102,11).Value = 917.6667
$ws.Cells.Item(102,12).Value = 2191
$ws.Cells.Item(102,13).Value = 704.3333
$ws.Cells.Item(102,14).Value = -5435

$ws.Cells.Item(132,8).Value = 1806
$ws.Cells.Item(132,9).Value = 1806
$ws.Cells.Item(132,10).Value = 0
$ws.Cells.Item(132,11).Value = 5418
$ws.Cells.Item(132,12).Value = 0
$ws.Cells.Item(132,13).Value = -2888

$ws.Cells.Item(136,8).Value = 3926.4
$ws.Cells.Item(136,9).Value = 3926.4
$ws.Cells.Item(136,10).Value = 0
$ws.Cells.Item(136,11).Value = 11779.2
$ws.Cells.Item(136,12).Value = 0
$ws.Cells.Item(136,13).Value = -9229.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80,8).Value = 293.5
$ws.Cells.Item(80,9).Value = 190
$ws.Cells.Item(80,10).Value = 314.2
$ws.Cells.Item(80,11).Value = 190
$ws.Cells.Item(80,12).Value = 314.2
$ws.Cells.Item(80,13).Value = 808
$ws.Cells.Item(80,14).Value = -2310.2

$ws.Cells.Item(83,8).Value = 293.5
$ws.Cells.Item(83,9).Value = 190
$ws.Cells.Item(83,10).Value = 314.2
$ws.Cells.Item(83,11).Value = 950
$ws.Cells.Item(83,12).Value = 1571
$ws.Cells.Item(83,13).Value = 4042
$ws.Cells.Item(83,14).Value = -11555

$ws.Cells.Item(86,8).Value = 4863.5713
$ws.Cells.Item(86,9).Value = 2886.25
$ws.Cells.Item(86,10).Value = 7500
$ws.Cells.Item(86,11).Value = 2886.25
$ws.Cells.Item(86,12).Value = 7500
$ws.Cells.Item(86,13).Value = -1763.25
$ws.Cells.Item(86,14).Value = -9746

$ws.Cells.Item(89,8).Value = 4863.5713
$ws.Cells.Item(89,9).Value = 2886.25
$ws.Cells.Item(89,10).Value = 7500
$ws.Cells.Item(89,11).Value = 14431.25
$ws.Cells.Item(89,12).Value = 37500
$ws.Cells.Item(89,13).Value = -8815.25
$ws.Cells.Item(89,14).Value = -48732

$ws.Cells.Item(105,8).Value = 3521.3333
$ws.Cells.Item(105,9).Value = 3613.2856
$ws.Cells.Item(105,10).Value = 3199.5
$ws.Cells.Item(105,11).Value = 3613.2856
$ws.Cells.Item(105,12).Value = 3199.5
$ws.Cells.Item(105,13).Value = -1866.2856

$ws.Cells.Item(134,8).Value = 8651.556
$ws.Cells.Item(134,9).Value = 8358.0625
$ws.Cells.Item(134,10).Value = 10999.5
$ws.Cells.Item(134,11).Value = 25074.1875
$ws.Cells.Item(134,12).Value = 32998.5
$ws.Cells.Item(134,13).Value = -22539.1875

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(36,8).Value = 149
$ws.Cells.Item(36,9).Value = 149
$ws.Cells.Item(36,10).Value = 0
$ws.Cells.Item(36,11).Value = 149
$ws.Cells.Item(36,12).Value = 0
$ws.Cells.Item(36,13).Value = 239

$ws.Cells.Item(40,8).Value = 149
$ws.Cells.Item(40,9).Value = 149
$ws.Cells.Item(40,10).Value = 0
$ws.Cells.Item(40,11).Value = 149
$ws.Cells.Item(40,12).Value = 0
$ws.Cells.Item(40,13).Value = 11

$ws.Cells.Item(62,8).Value = 3149
$ws.Cells.Item(62,9).Value = 0
$ws.Cells.Item(62,10).Value = 3149
$ws.Cells.Item(62,11).Value = 0
$ws.Cells.Item(62,12).Value = 3149
$ws.Cells.Item(62,13).ClearContents()
$ws.Cells.Item(62,14).Value = -4397

$ws.Cells.Item(65,8).Value = 3149
$ws.Cells.Item(65,9).Value = 0
$ws.Cells.Item(65,10).Value = 3149
$ws.Cells.Item(65,11).Value = 0
$ws.Cells.Item(65,12).Value = 15745
$ws.Cells.Item(65,13).ClearContents()
$ws.Cells.Item(65,14).Value = -21985

$ws.Cells.Item(105,8).Value = 1165
$ws.Cells.Item(105,9).Value = 820
$ws.Cells.Item(105,10).Value = 2200
$ws.Cells.Item(105,11).Value = 820
$ws.Cells.Item(105,12).Value = 2200
$ws.Cells.Item(105,13).Value = 927

$ws.Cells.Item(107,8).Value = 623.5
$ws.Cells.Item(107,9).Value = 249.5
$ws.Cells.Item(107,10).Value = 997.5
$ws.Cells.Item(107,11).Value = 249.5
$ws.Cells.Item(107,12).Value = 997.5
$ws.Cells.Item(107,13).Value = 1670.5
$ws.Cells.Item(107,14).Value = -4837.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2,8).Value = 230.33333
$ws.Cells.Item(2,9).Value = 48
$ws.Cells.Item(2,10).Value = 321.5
$ws.Cells.Item(2,11).Value = 288
$ws.Cells.Item(2,12).Value = 1929
$ws.Cells.Item(2,13).Value = -175
$ws.Cells.Item(2,14).Value = -2155

$ws.Cells.Item(12,8).Value = 101.375
$ws.Cells.Item(12,9).Value = 20.333334
$ws.Cells.Item(12,10).Value = 150
$ws.Cells.Item(12,11).Value = 61.000002
$ws.Cells.Item(12,12).Value = 450
$ws.Cells.Item(12,13).Value = 111.999998
$ws.Cells.Item(12,14).Value = -796

$ws.Cells.Item(17,8).Value = 841.25
$ws.Cells.Item(17,9).Value = 6
$ws.Cells.Item(17,10).Value = 1119.6666
$ws.Cells.Item(17,11).Value = 18
$ws.Cells.Item(17,12).Value = 3358.9998
$ws.Cells.Item(17,13).Value = 151
$ws.Cells.Item(17,14).Value = -3696.9998

$ws.Cells.Item(23,8).Value = 825.5
$ws.Cells.Item(23,9).Value = 591
$ws.Cells.Item(23,10).Value = 1060
$ws.Cells.Item(23,11).Value = 1773
$ws.Cells.Item(23,12).Value = 3180
$ws.Cells.Item(23,13).Value = -1538

$ws.Cells.Item(37,8).Value = 99702
$ws.Cells.Item(37,9).Value = 0
$ws.Cells.Item(37,10).Value = 99702
$ws.Cells.Item(37,11).Value = 0
$ws.Cells.Item(37,12).Value = 299106
$ws.Cells.Item(37,14).Value = -299330

$ws.Cells.Item(132,8).Value = 2388.2144
$ws.Cells.Item(132,9).Value = 1648.8462
$ws.Cells.Item(132,10).Value = 12000
$ws.Cells.Item(132,11).Value = 14839.6158
$ws.Cells.Item(132,12).Value = 108000
$ws.Cells.Item(132,13).Value = -12309.6158

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(46,8).Value = 18750
$ws.Cells.Item(46,9).Value = 7500
$ws.Cells.Item(46,10).Value = 20625
$ws.Cells.Item(46,11).Value = 7500
$ws.Cells.Item(46,12).Value = 20625
$ws.Cells.Item(46,13).Value = -7344
$ws.Cells.Item(46,14).Value = -20937

$ws.Cells.Item(57,8).Value = 18800
$ws.Cells.Item(57,9).Value = 4500
$ws.Cells.Item(57,10).Value = 24928.572
$ws.Cells.Item(57,11).Value = 4500
$ws.Cells.Item(57,12).Value = 24928.572
$ws.Cells.Item(57,13).Value = -3680
$ws.Cells.Item(57,14).Value = -26568.572

$ws.Cells.Item(70,8).Value = 0
$ws.Cells.Item(70,9).Value = 0
$ws.Cells.Item(70,10).Value = 0
$ws.Cells.Item(70,11).Value = 0
$ws.Cells.Item(70,12).Value = 0
$ws.Cells.Item(70,13).ClearContents()

$ws.Cells.Item(73,8).Value = 0
$ws.Cells.Item(73,9).Value = 0
$ws.Cells.Item(73,10).Value = 0
$ws.Cells.Item(73,11).Value = 0
$ws.Cells.Item(73,12).Value = 0
$ws.Cells.Item(73,13).ClearContents()

$ws.Cells.Item(80,8).Value = 2548.3333
$ws.Cells.Item(80,9).Value = 2548.3333
$ws.Cells.Item(80,10).Value = 0
$ws.Cells.Item(80,11).Value = 2548.3333
$ws.Cells.Item(80,12).Value = 0
$ws.Cells.Item(80,13).Value = -1550.3333

$ws.Cells.Item(83,8).Value = 2548.3333
$ws.Cells.Item(83,9).Value = 2548.3333
$ws.Cells.Item(83,10).Value = 0
$ws.Cells.Item(83,11).Value = 12741.6665
$ws.Cells.Item(83,12).Value = 0
$ws.Cells.Item(83,13).Value = -7749.666499999999

$ws.Cells.Item(113,8).Value = 4994
$ws.Cells.Item(113,9).Value = 4993
$ws.Cells.Item(113,10).Value = 4995
$ws.Cells.Item(113,11).Value = 4993
$ws.Cells.Item(113,12).Value = 4995
$ws.Cells.Item(113,13).Value = -2823

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16,8).Value = 699.5
$ws.Cells.Item(16,9).Value = 699.5
$ws.Cells.Item(16,10).Value = 0
$ws.Cells.Item(16,11).Value = 699.5
$ws.Cells.Item(16,12).Value = 0
$ws.Cells.Item(16,13).Value = -529.5

$ws.Cells.Item(46,8).Value = 3071.56
$ws.Cells.Item(46,9).Value = 1273.75
$ws.Cells.Item(46,10).Value = 3414
$ws.Cells.Item(46,11).Value = 1273.75
$ws.Cells.Item(46,12).Value = 3414
$ws.Cells.Item(46,13).Value = -1085.75
$ws.Cells.Item(46,14).Value = -3790

$ws.Cells.Item(132,8).Value = 14563.409
$ws.Cells.Item(132,9).Value = 10041.25
$ws.Cells.Item(132,10).Value = 19990
$ws.Cells.Item(132,11).Value = 30123.75
$ws.Cells.Item(132,12).Value = 59970
$ws.Cells.Item(132,13).Value = -27593.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4,8).Value = 2800.25
$ws.Cells.Item(4,9).Value = 1002
$ws.Cells.Item(4,10).Value = 3399.6667
$ws.Cells.Item(4,11).Value = 1002
$ws.Cells.Item(4,12).Value = 3399.6667
$ws.Cells.Item(4,13).Value = -889
$ws.Cells.Item(4,14).Value = -3625.6667

$ws.Cells.Item(100,8).Value = 1007.2222
$ws.Cells.Item(100,9).Value = 1138.5714
$ws.Cells.Item(100,10).Value = 547.5
$ws.Cells.Item(100,11).Value = 2277.1428
$ws.Cells.Item(100,12).Value = 1095
$ws.Cells.Item(100,13).Value = -1736.1428
